$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H7").Value = "this is dummy"
$ws.Range("H7").Select()
